$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 83, shifting existing rows 83:131 down to 84:132
$ws.Rows("83:83").Insert()

# Populate the newly inserted row 83 with the new weekly price record
$ws.Range("A83").Value = 11
$ws.Range("B83").Value = "Vega Monumental Concepción"
$ws.Range("C83").Value = "Bíobío"
$ws.Range("D83").Value = 44572
$ws.Range("E83").Value = 8
$ws.Range("F83").Value = "Fruta"
$ws.Range("G83").Value = 100108
$ws.Range("H83").Value = "Tropicales y subtropicales"
$ws.Range("I83").Value = 100108005
$ws.Range("J83").Value = "Piña"
$ws.Range("K83").Value = "Caramelo"
$ws.Range("L83").Value = "Primera"
$ws.Range("M83").Value = 200
$ws.Range("N83").Value = 15000
$ws.Range("O83").Value = 16000
$ws.Range("P83").Value = 15500
$ws.Range("Q83").Value = "$/caja 12 unidades"
$ws.Range("R83").Value = "Ecuador"
$ws.Range("S83").Value = 1292
$ws.Range("T83").Value = 12
